$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column B values for rows 1-4
$ws.Range("B1").Value = "John Tomanelli"
$ws.Range("B2").Value = "CFD89370C0"
$ws.Range("B3").Value = "IN -> 2017/02/15 14:46"
$ws.Range("B4").Value = "OUT -> 2017/02/15 14:47"

# Add new row 5 in column A
$ws.Range("A5").Value = "IN -> 2017/02/14 17:26"
